# Updated service entrance and calculations of installed and demand load
# for correct engine power.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Motores Elétricos" (Electric Motors) installed power corrected.
# Dependent formulas (D6 = B6*C6 and D7 = SUM(D2:D6)) recalc automatically.
$ws.Range("B6").Value = 18.41

# Selection moved from the old F12 leftover to the full data range.
$excel.Goto($ws.Range("A2:D7"))
